$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035143382025501
$ws.Cells.Item(2, 4).Value = 1.03813277289284
$ws.Cells.Item(2, 5).Value = 1.033982385555172
$ws.Cells.Item(2, 6).Value = 1.043601681888112
$ws.Cells.Item(2, 9).Value = 1.035898460803204
$ws.Cells.Item(2, 10).Value = 1.040258767711315
$ws.Cells.Item(2, 11).Value = 1.040921676020302
$ws.Cells.Item(2, 12).Value = 1.036783169240925
$ws.Cells.Item(2, 13).Value = 1.046375086769898

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.036899391494863
$ws.Cells.Item(3, 4).Value = 1.039492005783241
$ws.Cells.Item(3, 5).Value = 1.035522176796945
$ws.Cells.Item(3, 6).Value = 1.045435936546052
$ws.Cells.Item(3, 9).Value = 1.036417757860137
$ws.Cells.Item(3, 10).Value = 1.041654133013538
$ws.Cells.Item(3, 11).Value = 1.042089408047707
$ws.Cells.Item(3, 12).Value = 1.038130106514442
$ws.Cells.Item(3, 13).Value = 1.048017737731556

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.038031744876469
$ws.Cells.Item(4, 4).Value = 1.040368104538689
$ws.Cells.Item(4, 5).Value = 1.036514585823094
$ws.Cells.Item(4, 6).Value = 1.046619424254099
$ws.Cells.Item(4, 9).Value = 1.036750536247008
$ws.Cells.Item(4, 10).Value = 1.042552890941932
$ws.Cells.Item(4, 11).Value = 1.04284107270681
$ws.Cells.Item(4, 12).Value = 1.038997249882364
$ws.Cells.Item(4, 13).Value = 1.049076826862691

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.038506870827355
$ws.Cells.Item(5, 4).Value = 1.040735612094549
$ws.Cells.Item(5, 5).Value = 1.036930868047742
$ws.Cells.Item(5, 6).Value = 1.047116167797698
$ws.Cells.Item(5, 9).Value = 1.036889666701384
$ws.Cells.Item(5, 10).Value = 1.042929753695051
$ws.Cells.Item(5, 11).Value = 1.043156142837921
$ws.Cells.Item(5, 12).Value = 1.039360755311228
$ws.Cells.Item(5, 13).Value = 1.049521171726078

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.038586593426847
$ws.Cells.Item(6, 4).Value = 1.040797271448431
$ws.Cells.Item(6, 5).Value = 1.037000709807371
$ws.Cells.Item(6, 6).Value = 1.047199527163983
$ws.Cells.Item(6, 9).Value = 1.036912982374283
$ws.Cells.Item(6, 10).Value = 1.042992973870697
$ws.Cells.Item(6, 11).Value = 1.043208990345303
$ws.Cells.Item(6, 12).Value = 1.039421728832574
$ws.Cells.Item(6, 13).Value = 1.049595727129628

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.038038097104102
$ws.Cells.Item(7, 4).Value = 1.040373018334958
$ws.Cells.Item(7, 5).Value = 1.03652015182481
$ws.Cells.Item(7, 6).Value = 1.046626064862762
$ws.Cells.Item(7, 9).Value = 1.036752398330054
$ws.Cells.Item(7, 10).Value = 1.042557930411112
$ws.Cells.Item(7, 11).Value = 1.042845286325191
$ws.Cells.Item(7, 12).Value = 1.039002111131766
$ws.Cells.Item(7, 13).Value = 1.049082767718428

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.035737652975054
$ws.Cells.Item(8, 4).Value = 1.038592847383409
$ws.Cells.Item(8, 5).Value = 1.034503590730967
$ws.Cells.Item(8, 6).Value = 1.044222291686268
$ws.Cells.Item(8, 9).Value = 1.036074634761056
$ws.Cells.Item(8, 10).Value = 1.040731203718724
$ws.Cells.Item(8, 11).Value = 1.041317139139896
$ws.Cells.Item(8, 12).Value = 1.037239296460396
$ws.Cells.Item(8, 13).Value = 1.04693102893703

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031653242701107
$ws.Cells.Item(9, 4).Value = 1.035429192271198
$ws.Cells.Item(9, 5).Value = 1.030919242870372
$ws.Cells.Item(9, 6).Value = 1.039959668018876
$ws.Cells.Item(9, 9).Value = 1.0348552105366
$ws.Cells.Item(9, 10).Value = 1.037479906027776
$ws.Cells.Item(9, 11).Value = 1.038593633222114
$ws.Cells.Item(9, 12).Value = 1.034098484299864
$ws.Cells.Item(9, 13).Value = 1.04310938272528

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028908446275615
$ws.Cells.Item(10, 4).Value = 1.033301239726247
$ws.Cells.Item(10, 5).Value = 1.028507851288905
$ws.Cells.Item(10, 6).Value = 1.037098664444674
$ws.Cells.Item(10, 9).Value = 1.034024961030668
$ws.Cells.Item(10, 10).Value = 1.035289638786459
$ws.Cells.Item(10, 11).Value = 1.036756508445603
$ws.Cells.Item(10, 12).Value = 1.031980433748075
$ws.Cells.Item(10, 13).Value = 1.040540342207647

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027714467743331
$ws.Cells.Item(11, 4).Value = 1.032375150395715
$ws.Cells.Item(11, 5).Value = 1.027458280989575
$ws.Cells.Item(11, 6).Value = 1.035854984633391
$ws.Cells.Item(11, 9).Value = 1.033661257734151
$ws.Cells.Item(11, 10).Value = 1.034335618579828
$ws.Cells.Item(11, 11).Value = 1.035955744095018
$ws.Cells.Item(11, 12).Value = 1.03105734166909
$ws.Cells.Item(11, 13).Value = 1.039422625868098

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027270126935031
$ws.Cells.Item(12, 4).Value = 1.032030441038021
$ws.Cells.Item(12, 5).Value = 1.027067588962826
$ws.Cells.Item(12, 6).Value = 1.03539227533561
$ws.Cells.Item(12, 9).Value = 1.033525523069557
$ws.Cells.Item(12, 10).Value = 1.033980389742597
$ws.Cells.Item(12, 11).Value = 1.035657496034356
$ws.Cells.Item(12, 12).Value = 1.030713549438907
$ws.Cells.Item(12, 13).Value = 1.039006637836862

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027365478115905
$ws.Cells.Item(13, 4).Value = 1.032104415151007
$ws.Cells.Item(13, 5).Value = 1.027151431838194
$ws.Cells.Item(13, 6).Value = 1.035491562420556
$ws.Cells.Item(13, 9).Value = 1.033554667666212
$ws.Cells.Item(13, 10).Value = 1.034056626904434
$ws.Cells.Item(13, 11).Value = 1.035721508097531
$ws.Cells.Item(13, 12).Value = 1.030787335740193
$ws.Cells.Item(13, 13).Value = 1.039095906063794

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027677755746472
$ws.Cells.Item(14, 4).Value = 1.032346671367995
$ws.Cells.Item(14, 5).Value = 1.027426003414685
$ws.Cells.Item(14, 6).Value = 1.0358167523985
$ws.Cells.Item(14, 9).Value = 1.033650050948459
$ws.Cells.Item(14, 10).Value = 1.034306272973776
$ws.Cells.Item(14, 11).Value = 1.035931107402533
$ws.Cells.Item(14, 12).Value = 1.031028942467032
$ws.Cells.Item(14, 13).Value = 1.039388256988105

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027870047761474
$ws.Cells.Item(15, 4).Value = 1.032495837702086
$ws.Cells.Item(15, 5).Value = 1.02759506478206
$ws.Cells.Item(15, 6).Value = 1.03601701243253
$ws.Cells.Item(15, 9).Value = 1.033708734809938
$ws.Cells.Item(15, 10).Value = 1.034459973196712
$ws.Cells.Item(15, 11).Value = 1.036060140841227
$ws.Cells.Item(15, 12).Value = 1.031177682562597
$ws.Cells.Item(15, 13).Value = 1.039568274947282

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028987569616031
$ws.Cells.Item(16, 4).Value = 1.033362601378587
$ws.Cells.Item(16, 5).Value = 1.028577391853307
$ws.Cells.Item(16, 6).Value = 1.037181099303228
$ws.Cells.Item(16, 9).Value = 1.034049009647656
$ws.Cells.Item(16, 10).Value = 1.035352833949889
$ws.Cells.Item(16, 11).Value = 1.036809540019816
$ws.Cells.Item(16, 12).Value = 1.032041569098322
$ws.Cells.Item(16, 13).Value = 1.0406144077711

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029687082532651
$ws.Cells.Item(17, 4).Value = 1.033905036865194
$ws.Cells.Item(17, 5).Value = 1.029192113836582
$ws.Cells.Item(17, 6).Value = 1.037909986499785
$ws.Cells.Item(17, 9).Value = 1.034261324983688
$ws.Cells.Item(17, 10).Value = 1.035911384000781
$ws.Cells.Item(17, 11).Value = 1.037278194464189
$ws.Cells.Item(17, 12).Value = 1.032581852728083
$ws.Cells.Item(17, 13).Value = 1.041269183963483

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.030094571008382
$ws.Cells.Item(18, 4).Value = 1.034220980679787
$ws.Cells.Item(18, 5).Value = 1.02955014902686
$ws.Cells.Item(18, 6).Value = 1.038334668041716
$ws.Cells.Item(18, 9).Value = 1.034384760262101
$ws.Cells.Item(18, 10).Value = 1.036236635833551
$ws.Cells.Item(18, 11).Value = 1.037551044570246
$ws.Cells.Item(18, 12).Value = 1.032896417055786
$ws.Cells.Item(18, 13).Value = 1.041650593459926

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.030233425598498
$ws.Cells.Item(19, 4).Value = 1.034328633591618
$ws.Cells.Item(19, 5).Value = 1.029672141888438
$ws.Cells.Item(19, 6).Value = 1.038479395121429
$ws.Cells.Item(19, 9).Value = 1.034426780107138
$ws.Cells.Item(19, 10).Value = 1.036347447114954
$ws.Cells.Item(19, 11).Value = 1.037643993693671
$ws.Cells.Item(19, 12).Value = 1.033003578577145
$ws.Cells.Item(19, 13).Value = 1.041780558316519

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.029612085946403
$ws.Cells.Item(20, 4).Value = 1.033846885298143
$ws.Cells.Item(20, 5).Value = 1.029126214098114
$ws.Cells.Item(20, 6).Value = 1.037831832171824
$ws.Cells.Item(20, 9).Value = 1.034238587472201
$ws.Cells.Item(20, 10).Value = 1.035851512956866
$ws.Cells.Item(20, 11).Value = 1.037227964946054
$ws.Cells.Item(20, 12).Value = 1.032523944913236
$ws.Cells.Item(20, 13).Value = 1.041198985620009

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027585821254372
$ws.Cells.Item(21, 4).Value = 1.032275352920486
$ws.Cells.Item(21, 5).Value = 1.027345172140363
$ws.Cells.Item(21, 6).Value = 1.035721012980361
$ws.Cells.Item(21, 9).Value = 1.033621980656172
$ws.Cells.Item(21, 10).Value = 1.034232782411808
$ws.Cells.Item(21, 11).Value = 1.035869408077713
$ws.Cells.Item(21, 12).Value = 1.030957820712591
$ws.Cells.Item(21, 13).Value = 1.039302189693357

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026306933272757
$ws.Cells.Item(22, 4).Value = 1.031283102432816
$ws.Cells.Item(22, 5).Value = 1.026220519720392
$ws.Cells.Item(22, 6).Value = 1.034389499042844
$ws.Cells.Item(22, 9).Value = 1.033230593940126
$ws.Cells.Item(22, 10).Value = 1.033210018190491
$ws.Cells.Item(22, 11).Value = 1.035010543635079
$ws.Cells.Item(22, 12).Value = 1.029967834329487
$ws.Cells.Item(22, 13).Value = 1.038104852879946

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.026985367794485
$ws.Cells.Item(23, 4).Value = 1.031809513673017
$ws.Cells.Item(23, 5).Value = 1.026817184939047
$ws.Cells.Item(23, 6).Value = 1.035095780677984
$ws.Cells.Item(23, 9).Value = 1.03343842904051
$ws.Cells.Item(23, 10).Value = 1.033752685564648
$ws.Cells.Item(23, 11).Value = 1.035466293382779
$ws.Cells.Item(23, 12).Value = 1.030493153806663
$ws.Cells.Item(23, 13).Value = 1.038740040867493

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029645975298942
$ws.Cells.Item(24, 4).Value = 1.033873162877171
$ws.Cells.Item(24, 5).Value = 1.029155992964969
$ws.Cells.Item(24, 6).Value = 1.037867148187955
$ws.Cells.Item(24, 9).Value = 1.034248862825131
$ws.Cells.Item(24, 10).Value = 1.035878567775574
$ws.Cells.Item(24, 11).Value = 1.037250663072848
$ws.Cells.Item(24, 12).Value = 1.032550112735215
$ws.Cells.Item(24, 13).Value = 1.04123070680629

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.03271292461277
$ws.Cells.Item(25, 4).Value = 1.036250332009258
$ws.Cells.Item(25, 5).Value = 1.031849652369432
$ws.Cells.Item(25, 6).Value = 1.041064963985422
$ws.Cells.Item(25, 9).Value = 1.035173479614525
$ws.Cells.Item(25, 10).Value = 1.038324379774167
$ws.Cells.Item(25, 11).Value = 1.039301446699599
$ws.Cells.Item(25, 12).Value = 1.03491464844589
$ws.Cells.Item(25, 13).Value = 1.044101038395631
